$wb = $excel.ActiveWorkbook

# The CMS sheet (last sheet) gets a new data row, duplicating row 2 but
# with a different "Contact_ID" value (A column).
$ws = $wb.Worksheets.Item("CMS")

$ws.Rows(2).Copy()
$ws.Rows(4).Insert(-4121)
$excel.CutCopyMode = $false

# Row-insert doesn't carry the (format-only) wrap-text style that column A
# uses, so reapply it before writing the real value.
$ws.Range("A4").WrapText = $true
$ws.Range("A4").Value = 456

# Select rows 5:6 (entire rows) with A5 as the active cell, and make the
# CMS sheet the active / selected tab.
$ws.Range("A5:XFD6").Select()
$ws.Activate()

$wb.Save()
